$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Mid-term date moved from Friday, July 13th, 2022 to Friday, July 15th, 2022
$ws.Range("B16").Value = "Friday, July 15th, 2022"
